# Apply updated cryptocurrency price/volume data to Sheet1.
# Column D ("Price") values are prefixed with a leading apostrophe so Excel
# stores them as text (matching the source data) instead of auto-converting
# decimal-looking strings into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '''65.486.34'
$ws.Cells.Item(2, 5).Value = '  -1.79%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '''3.502.46'
$ws.Cells.Item(3, 5).Value = '  -2.55%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.03%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''598.01'
$ws.Cells.Item(5, 5).Value = '  -2.02%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''141.72'
$ws.Cells.Item(6, 5).Value = '  -4.61%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''3.500.72'
$ws.Cells.Item(7, 5).Value = '  -2.51%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.999'
$ws.Cells.Item(8, 5).Value = '  -0.27%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.509'
$ws.Cells.Item(9, 5).Value = '  +3.75%  '

# Row 10
$ws.Cells.Item(10, 4).Value = '''7.77'
$ws.Cells.Item(10, 5).Value = '  -2.83%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.129'
$ws.Cells.Item(11, 5).Value = '  -5.49%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''0.400'
$ws.Cells.Item(12, 5).Value = '  -3.77%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''4.096.57'
$ws.Cells.Item(13, 5).Value = '  -2.59%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '''0.0000191'
$ws.Cells.Item(14, 5).Value = '  -8.43%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '''28.10'
$ws.Cells.Item(15, 5).Value = '  -6.76%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''3.501.85'
$ws.Cells.Item(16, 5).Value = '  -2.10%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +1.21%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '''65.437.93'
$ws.Cells.Item(18, 5).Value = '  -2.00%  '

# Row 19
$ws.Cells.Item(19, 4).Value = '''10.76'
$ws.Cells.Item(19, 5).Value = '  -6.11%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -3.59%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''14.33'
$ws.Cells.Item(21, 5).Value = '  -5.20%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''414.26'
$ws.Cells.Item(22, 5).Value = '  -4.09%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  -5.94%  '

# Row 24
$ws.Cells.Item(24, 4).Value = '''76.57'
$ws.Cells.Item(24, 5).Value = '  -3.31%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''3.643.38'
$ws.Cells.Item(25, 5).Value = '  -2.46%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.03%  '

# Row 27
$ws.Cells.Item(27, 5).Value = '  -7.03%  '

# Row 28
$ws.Cells.Item(28, 5).Value = '  -3.45%  '

# Row 29
$ws.Cells.Item(29, 4).Value = '''7.67'
$ws.Cells.Item(29, 5).Value = '  -6.85%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''8.77'
$ws.Cells.Item(30, 5).Value = '  -6.01%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  -0.03%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '''3.509.29'
$ws.Cells.Item(32, 5).Value = '  -2.31%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -1.92%  '

# Row 34
$ws.Cells.Item(34, 5).Value = '  -5.93%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -0.02%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  -8.26%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''7.44'
$ws.Cells.Item(37, 5).Value = '  -5.35%  '

# Row 38
$ws.Cells.Item(38, 4).Value = '''175.41'
$ws.Cells.Item(38, 5).Value = '  +0.42%  '

# Row 39 (ImmutableX)
$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(39, 4).Value = '''1.57'
$ws.Cells.Item(39, 5).Value = '  -8.50%  '

# Row 40 (NEARProtocol)
$ws.Cells.Item(40, 2).Value = 'NEARProtocol'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(40, 4).Value = '''5.15'
$ws.Cells.Item(40, 5).Value = '  -8.87%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''0.0805'
$ws.Cells.Item(41, 5).Value = '  -6.28%  '

# Row 42 (Filecoin)
$ws.Cells.Item(42, 2).Value = 'Filecoin'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(42, 4).Value = '''4.94'
$ws.Cells.Item(42, 5).Value = '  -5.86%  '

# Row 43 (Mantle)
$ws.Cells.Item(43, 2).Value = 'Mantle'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(43, 4).Value = '''0.851'
$ws.Cells.Item(43, 5).Value = '  -5.06%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''45.16'
$ws.Cells.Item(44, 5).Value = '  -2.12%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''1.74'

# Row 46
$ws.Cells.Item(46, 5).Value = '  +0.06%  '

# Row 47 (dogwifhat)
$ws.Cells.Item(47, 2).Value = 'dogwifhat'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(47, 4).Value = '''2.31'
$ws.Cells.Item(47, 5).Value = '  -9.75%  '

# Row 48 (Cosmos)
$ws.Cells.Item(48, 2).Value = 'Cosmos'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(48, 4).Value = '''6.97'
$ws.Cells.Item(48, 5).Value = '  -3.59%  '

# Row 49 (EnergySwap)
$ws.Cells.Item(49, 2).Value = 'EnergySwap'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Cells.Item(49, 4).Value = '''22.91'
$ws.Cells.Item(49, 5).Value = '  -4.67%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  -9.58%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '''0.896'
$ws.Cells.Item(51, 5).Value = '  -5.98%  '
